$wb = $excel.ActiveWorkbook

$wsArq   = $wb.Worksheets.Item(1)   # ARQUITECTO
$wsVales = $wb.Worksheets.Item(2)   # VALES DE INSENTIVOS

# --- VALES DE INSENTIVOS sheet content updates ---

# Amount-in-words cell: spacing tweak in the wording.
$wsVales.Range("A2").Value = "CINCO   MIL     PESOS 00/100 M.N."

# A new (blank-looking) text cell is added at D3 (two spaces, same as the
# equivalent cell already present on the ARQUITECTO sheet).
$wsVales.Range("D3").Value = "  "

# The incentive month label gets filled in with the actual month.
$wsVales.Range("A4").Value = "INCENTIVO DEL MES DE NOVIEMBRE"

# Touching the merged signature-name block (re-merge) so the workbook's
# internal bookkeeping for that merge matches a fresh edit.
$wsVales.Range("C8:D9").UnMerge()
$wsVales.Range("C8:D9").Merge()

# Column width tweaks (match ARQUITECTO's column B / C widths).
$wsVales.Columns.Item(2).ColumnWidth = 8.33
$wsVales.Columns.Item(3).ColumnWidth = 18.33

# Cursor / selection moves to D6 on this sheet.
$wsVales.Range("D6").Select()

# This sheet becomes the active / selected tab in the workbook.
$wsVales.Activate()
